# Update Work Week and Social Spending
# Rewrites the GDP-per-Capita "Data" sheet: replaces the historical
# 1973-2010 values with a revised series, and extends the table through
# 2016 with newly published figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Revised data: year -> GDP per Capita value (string form preserved
# exactly as published, numeric where the source gives an integer).
$values = [ordered]@{
    1973 = "6527"
    1974 = ""
    1975 = ""
    1976 = ""
    1977 = ""
    1978 = ""
    1979 = ""
    1980 = "5746"
    1981 = "5622"
    1982 = "5498"
    1983 = "5547"
    1984 = "5469"
    1985 = "5364"
    1986 = "5319"
    1987 = "5040"
    1988 = "5515"
    1989 = "4962"
    1990 = "4803"
    1991 = "4270.37012179731"
    1992 = "2832.70342443527"
    1993 = "2482.41359224165"
    1994 = "1928.3185195337"
    1995 = "1666.18395821173"
    1996 = "1559.52118921367"
    1997 = "1552.63236017053"
    1998 = "1589.33155511695"
    1999 = "1621.3205712536"
    2000 = "1723.35632964691"
    2001 = "1845.47821894097"
    2002 = "1982.88366574433"
    2003 = "2135.90828146068"
    2004 = "2301.29208623307"
    2005 = "2416.86704371227"
    2006 = "2524.84473749478"
    2007 = "2647.20938715616"
    2008 = "2801.77617916365"
    2009 = "2853.88848045587"
    2010 = "2975.92680864722"
    2011 = "3111"
    2012 = "3285"
    2013 = "3471"
    2014 = "3633"
    2015 = "3789"
    2016 = "3991"
}

foreach ($year in $values.Keys) {
    $row = $year - 1973 + 2

    $ws.Cells.Item($row, 1).Value = 762
    $ws.Cells.Item($row, 2).Value = "Tajikistan"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $year

    $newValue = $values[$year]
    if ($newValue -eq "") {
        # These rows were already blank in the workbook (no published
        # figure for that year) and stay blank after the update, so
        # there is nothing to write here.
        continue
    }

    # The Data column stores its numbers as plain text (shared strings),
    # matching the published source file. Force text storage by
    # pre-formatting the cell, then strip the formatting back off so no
    # stray number-format style is left behind on the cell.
    $cell = $ws.Cells.Item($row, 5)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.ClearFormats()
}
